$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BMU_Inmates")

# --- Row 2: Maier, Christopher -> Dingle, Derp ---
$ws.Range("A2").Value = "Dingle, Derp"
$ws.Range("B2").Value = 12345
$ws.Range("D2").Value = 8

# --- Row 3: Stubbs, Brian -> Tweedle, Dee ---
$ws.Range("A3").Value = "Tweedle, Dee"
$ws.Range("B3").Value = 12346

# --- Row 4: Riley, Cordero -> Ship, Lee ---
$ws.Range("A4").Value = "Ship, Lee"
$ws.Range("B4").Value = 12347
$ws.Range("D4").Value = 1

# --- Row 5: previously blank, now populated with a new entry ---
# Copy the date cell formatting from E2 so E5 picks up the same
# number format / style (rather than Excel's auto-detected date style).
$ws.Range("E2").Copy()
$ws.Range("E5").PasteSpecial(-4122)

$ws.Range("A5").Value = "Shmo, Joe"
$ws.Range("B5").Value = 12348
$ws.Range("D5").Value = 12
$ws.Range("E5").Value = "11/16/2021"
$ws.Range("F5").Value = "M. Ham"

# --- Update the active selection shown in the sheet view ---
$ws.Activate() | Out-Null
$ws.Range("F5").Select() | Out-Null
